$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo "tshirt" -> "t-shirt" in the "Catégorie" definition cell (B3),
# which is part of the glossary table (Tableau2).
$cell = $ws.Range("B3")
$cell.Value = "Catégorie à laquelle appartient un vêtement : pantalon, jupe, veste, t-shirt, blouse… - un vêtement appartient à une seule catégorie. L'utilisateur peut ajouter des catégories."

# Leave the final selection where the author ended up (E7), matching the
# cursor position recorded in the saved workbook.
$ws.Range("E7").Select()

$wb.Save()
